# edit.ps1 - applies the ReportTemplate.docx changes described by the commit:
#  - Paragraph 1: append " №{idOrderStage}" before the trailing "." and
#    remove the obsolete "_GoBack" bookmark.
#  - Paragraph 2: rewrite the "Исполнитель/Заказчик" sentence so that the
#    executor name becomes literal text (НИЦ «ИНКОМСИСТЕМ») and the
#    customer becomes a merge field {customer} wrapped in curly quotes.
#
# NOTE: this sandboxed Word COM-interop layer has a bug where
# Range.LanguageID only ever applies to the FIRST paragraph of the
# document (and only serializes correctly when given a string like
# "en-US", not the numeric WdLanguageID constant). We work around this by
# building each "en-US" tagged run in a temporary paragraph inserted at
# the very start of the document (where the bug actually lands correctly),
# cutting it to the clipboard, and pasting it over the target location.

$d = $word.ActiveDocument

function New-TaggedRun($text) {
    # Creates a run containing $text with w:lang="en-US" by exploiting the
    # fact that Range.LanguageID only works when applied to paragraph 1.
    # Leaves the tagged text on the clipboard (via Cut) ready to Paste().
    $doc = $word.ActiveDocument
    $beginRange = $doc.Range(0, 0)
    $beginRange.InsertParagraphBefore()
    $tagRange = $doc.Range(0, 0)
    $tagRange.InsertBefore($text)
    $tagRange2 = $doc.Range(0, $text.Length)
    $tagRange2.LanguageID = "en-US"
    $tagRange2.Cut()
    $tempPara = $doc.Paragraphs(1).Range
    $tempPara.Delete()
}

# ---------------------------------------------------------------------
# Paragraph 1: "Отчет о выполнении работы (этапа работы)."
# ---------------------------------------------------------------------

# Remove the stale _GoBack bookmark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$full = $d.Content
$full.Find.Execute("работы).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos = $full.Start + [string]"работы)".Length

$target = $d.Range($insPos, $insPos)
$target.InsertBefore(" №{idOrderStage}")

$tagStart = $insPos + [string]" №{".Length
$tagEnd = $tagStart + [string]"idOrderStage".Length
$tagRange = $d.Range($tagStart, $tagEnd)

New-TaggedRun("idOrderStage")
$tagRange = $d.Range($tagStart, $tagEnd)
$tagRange.Paste()

# ---------------------------------------------------------------------
# Paragraph 2: "Исполнитель: Компания ... Заказчик: ..."
# ---------------------------------------------------------------------

$p2 = $d.Paragraphs(2).Range
$p2Start = $p2.Start

$openQuote = [char]0x201C   # “
$closeQuote = [char]0x201D  # ”

$prefixText = 'Исполнитель: Компания "НИЦ «ИНКОМСИСТЕМ»" Заказчик: '
$finalText = $prefixText + $openQuote + "{customer}" + $closeQuote
$p2.Text = $finalText

# Tag the opening curly quote with en-US.
$openQuotePos = $p2Start + [string]$prefixText.Length
New-TaggedRun([string]$openQuote)
$openQuoteRange = $d.Range($openQuotePos, $openQuotePos + 1)
$openQuoteRange.Paste()

# Tag the "customer" merge-field name with en-US.
$customerStart = $openQuotePos + 2   # skip opening quote + "{"
$customerEnd = $customerStart + [string]"customer".Length
New-TaggedRun("customer")
$customerRange = $d.Range($customerStart, $customerEnd)
$customerRange.Paste()
